# Update automatico via Actualizar 02-06-2021 04-39-05
#
# Appends one more "availability check" block (14 rows) to Sheet1, mirroring
# the existing repeating 14-row pattern (Odoo, Blackbox, PowerBI, Dropbox,
# Odoo, GEE, UtilidadesOdoo, Filtros Dashboard, MapStore, GeoServer, Tomcat,
# Shiny, Github, EZ Exporter) with a new timestamp, and wires up the matching
# hyperlinks in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous block (rows 996-1009) gets its "Fecha" timestamp nudged by a
# couple of microseconds too (re-check landed a few ms later than the write).
for ($row = 996; $row -le 1009; $row++) {
    $ws.Cells.Item($row, 4).Value = 44233.1725618287
}

$timestamp = 44233.19375778658

# Name (col A), URL (col B) for each of the 14 rows in one cycle.
$rowsData = @(
    @{ Name = "Odoo";              Url = "https://www.dataintelligence-group.com/" },
    @{ Name = "Blackbox";          Url = "https://serviciodashboard.azurewebsites.net/" },
    @{ Name = "PowerBI";           Url = "https://powerbi.microsoft.com/es-es/" },
    @{ Name = "Dropbox";           Url = "https://www.dropbox.com/" },
    @{ Name = "Odoo";              Url = "https://dataintelligence.store/" },
    @{ Name = "GEE";               Url = "https://app-data-i.users.earthengine.app/" },
    @{ Name = "UtilidadesOdoo";    Url = "https://odooutil.azurewebsites.net/" },
    @{ Name = "Filtros Dashboard"; Url = "https://filtradordashboard.azurewebsites.net/" },
    @{ Name = "MapStore";          Url = "https://ide.dataintelligence-group.com/mapstore/#/" },
    @{ Name = "GeoServer";         Url = "https://ide.dataintelligence-group.com/geoserver/web/?0" },
    @{ Name = "Tomcat";            Url = "https://ide.dataintelligence-group.com/" },
    @{ Name = "Shiny";             Url = "https://rpubs.com/dataintelligence/" },
    @{ Name = "Github";            Url = "https://github.com/Sud-Austral/" },
    @{ Name = "EZ Exporter";       Url = "https://ezexporter.highviewapps.com/exports/export-profile/" }
)

$startRow = 1010

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $row = $startRow + $i
    $entry = $rowsData[$i]

    $nameCell = $ws.Cells.Item($row, 1)
    $urlCell = $ws.Cells.Item($row, 2)
    $availCell = $ws.Cells.Item($row, 3)
    $dateCell = $ws.Cells.Item($row, 4)

    $nameCell.Value = $entry.Name
    $availCell.Value = "Disponible"
    $dateCell.Value = $timestamp
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    if ($entry.Url -like "*#/") {
        # The MapStore link carries a '#/' fragment: Excel stores that as the
        # hyperlink's SubAddress (location="/") with the Target trimmed.
        $base = $entry.Url.Substring(0, $entry.Url.Length - 2)
        $urlCell.Value = $entry.Url
        $ws.Hyperlinks.Add($urlCell, $base, "/")
    } else {
        $urlCell.Value = $entry.Url
        $ws.Hyperlinks.Add($urlCell, $entry.Url)
    }

    # Hyperlinks.Add() re-styles the cell with a freshly minted xf record;
    # put it back on the shared "Hyperlink" cell style used by every other
    # link cell in column B.
    $urlCell.Style = "Hyperlink"
}
